$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Fourier-series coefficient data (dynamics identification re-run)
$ws.Range("A1").Value = 0.00333938462730028
$ws.Range("A2").Value = 0.00333384923609884
$ws.Range("A3").Value = -0.0483427458717823
$ws.Range("A4").Value = -0.260456087909828
$ws.Range("A5").Value = 0.300859661366
$ws.Range("A6").Value = 0.0489348920394959
$ws.Range("A7").Value = 0.029475024446471
$ws.Range("A8").Value = -0.153092694500051
$ws.Range("A9").Value = -0.634042080077492
$ws.Range("A10").Value = 0.576915794297815
$ws.Range("A11").Value = -0.0982900850264507
$ws.Range("A12").Value = 0.0168622783694167
$ws.Range("A13").Value = -0.294824039207678
$ws.Range("A14").Value = -0.256542767097894
$ws.Range("A15").Value = 0.294017335906422
$ws.Range("A16").Value = 0.240015488168106
$ws.Range("A17").Value = 0.121912319763234
$ws.Range("A18").Value = -0.163786511095751
$ws.Range("A19").Value = 0.0298695983684544
$ws.Range("A20").Value = 0.113952178280737
$ws.Range("A21").Value = -0.0677709101145603
$ws.Range("A22").Value = -1.36269842684397
$ws.Range("A23").Value = 0.00720203280867843
$ws.Range("A24").Value = 0.106061171029263
$ws.Range("A25").Value = 0.572283564124571
$ws.Range("A26").Value = -0.0228043935037851
$ws.Range("A27").Value = -0.662615245731844
$ws.Range("A28").Value = -0.138223596050814
$ws.Range("A29").Value = 0.170728474434344
$ws.Range("A30").Value = -0.0205234267280683
$ws.Range("A31").Value = 0.0592439574360344
$ws.Range("A32").Value = -0.0757318073806069
$ws.Range("A33").Value = -0.192509243285872
$ws.Range("A34").Value = -0.00856878590228507
$ws.Range("A35").Value = -0.00335629421967389
$ws.Range("A36").Value = 0.0049403279019122
$ws.Range("A37").Value = 0.00477108920676527
$ws.Range("A38").Value = 0.00249855438070927
$ws.Range("A39").Value = 0.0141936891896592
$ws.Range("A40").Value = 0.0606428150431505
$ws.Range("A41").Value = -0.189228986429901
$ws.Range("A42").Value = -0.535933094740337
$ws.Range("A43").Value = 0.51501576930382
$ws.Range("A44").Value = -1.72841212132103
$ws.Range("A45").Value = 0.0309928744032121
$ws.Range("A46").Value = 0.0231946587797607
$ws.Range("A47").Value = 0.532433526864129
$ws.Range("A48").Value = 0.0993388500174188
$ws.Range("A49").Value = -0.685814933112209
$ws.Range("A50").Value = 0.0472789852553486
$ws.Range("A51").Value = -0.158064460956765
$ws.Range("A52").Value = 0.33475202890947
$ws.Range("A53").Value = -0.175539468648756
$ws.Range("A54").Value = -0.00672581926646901
$ws.Range("A55").Value = 0.110485505289996
$ws.Range("A56").Value = 0.0509580994017887
$ws.Range("A57").Value = -0.191448426500055
$ws.Range("A58").Value = -0.244909277698867
$ws.Range("A59").Value = 0.215353666982259
$ws.Range("A60").Value = 0.170094141643033
$ws.Range("A61").Value = 0.458711201908083
$ws.Range("A62").Value = -0.738979774779363
$ws.Range("A63").Value = -0.54694756176714
$ws.Range("A64").Value = 0.378313949612719
$ws.Range("A65").Value = 0.229283012855533
$ws.Range("A66").Value = 0.152315742160172

# Select the full data range, matching the author's "select all" before saving
$ws.UsedRange.Select()
